$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.865412354469299
$ws.Range("B1").Value = 1.837787270545959
$ws.Range("C1").Value = 2.173841238021851
$ws.Range("D1").Value = 2.382610321044922
$ws.Range("E1").Value = 1.663984298706055
